$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.834.08'
$ws.Range('E2').Value = '  -3.97%  '
$ws.Range('D3').Value = '2.424.90'
$ws.Range('E3').Value = '  -7.05%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '538.83'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -6.08%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.32'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -7.42%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.580'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -7.49%  '
$ws.Range('D9').Value = '2.427.58'
$ws.Range('E9').Value = '  -6.89%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.104'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -10.19%  '
$ws.Range('E11').Value = '  -1.97%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.35'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -6.86%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.346'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -8.83%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.56'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -9.84%  '
$ws.Range('D15').Value = '2.865.02'
$ws.Range('E15').Value = '  -7.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000162'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -10.71%  '
$ws.Range('D17').Value = '60.777.86'
$ws.Range('E17').Value = '  -3.97%  '
$ws.Range('D18').Value = '2.428.24'
$ws.Range('E18').Value = '  -7.18%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.84'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -9.71%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.86'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -9.16%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.11'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -7.71%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '315.06'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -7.26%  '
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.82'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.78%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '63.05'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Value = '2.556.97'
$ws.Range('E26').Value = '  -6.96%  '
$ws.Range('E27').Value = '  -13.75%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.10%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.10'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -11.09%  '
$ws.Range('B30').Value = 'Fetch.AI'
$ws.Range('C30').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.43'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -12.00%  '
$ws.Range('B31').Value = 'Aptos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.56'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -6.94%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '513.79'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -8.59%  '
$ws.Range('E33').Value = '  -8.23%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.86'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -8.28%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.55'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -8.98%  '
$ws.Range('E36').Value = '  -0.13%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.59'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -13.88%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.70'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -11.70%  '
$ws.Range('E39').Value = '  -6.47%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.18'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.74%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '141.75'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -8.04%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.73'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -7.32%  '
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '40.14'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.18%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.21'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -11.35%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '138.92'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -13.47%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.56'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -8.11%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '20.79'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -12.59%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0526'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -9.30%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.579'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -7.32%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0924'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -6.71%  '
